# Update attendance_percentage, total_classes and classes_attended
# values for student_details (dynamic rendering update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# roll_no 2205460 - Bhaskar Lalwani
$ws.Range("C2").Value = 86
$ws.Range("D2").Value = 48
$ws.Range("E2").Value = 41

# roll_no 2205568 - Mayur Gogoi
$ws.Range("C3").Value = 65.40000000000001
$ws.Range("D3").Value = 46
$ws.Range("E3").Value = 30

# roll_no 2205533 - Aniruddha Mukherjee
$ws.Range("C4").Value = 70.7
$ws.Range("D4").Value = 40
$ws.Range("E4").Value = 28

# roll_no 2205967 - Amandeep Chourasia
$ws.Range("C5").Value = 79.59999999999999
$ws.Range("D5").Value = 46
$ws.Range("E5").Value = 36

# roll_no 2205557 - Ishaan Mukherjee
$ws.Range("C6").Value = 92.59999999999999
$ws.Range("D6").Value = 49
$ws.Range("E6").Value = 45
